# Input grid centroids for agricultural samples
# Update habitat/area metadata for the set of rows whose mfd_hab1/mfd_hab2
# were mis-coded as marine ("Saltwater" / "Open sea and tidal areas")
# when they should be freshwater lake samples. Also corrects the
# habitat_typenumber (F) from 1150 -> 3170 and adds the missing mfd_hab3
# (column P) value "Lake shores and lagoons".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(13, 14, 24, 26, 76, 82, 96, 98, 101, 119, 127, 156, 163, 173, 218, 228, 261)

foreach ($r in $rows) {
    # habitat_typenumber: 1150 -> 3170 (must stay text, matching the rest of
    # the column, so force it via a leading apostrophe then strip the
    # resulting "quote prefix" style back to Normal so no stray formatting
    # is introduced).
    $ws.Range("F$r").Value = "'3170"
    $ws.Range("F$r").Style = "Normal"

    # mfd_hab1: Saltwater -> Freshwater
    $ws.Range("N$r").Value = "Freshwater"

    # mfd_hab2: Open sea and tidal areas -> Standing freshwater
    $ws.Range("O$r").Value = "Standing freshwater"

    # mfd_hab3: newly populated with Lake shores and lagoons
    $ws.Range("P$r").Value = "Lake shores and lagoons"
}
